# Agrega la fila de resultados de Alyne Corona (SmartScore) al final de la hoja
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Formula = "Alyne Corona_20251128_162856"
$ws.Range("B18").Formula = "'"
$ws.Range("C18").Formula = "Alyne Corona"
$ws.Range("D18").Formula = 1
$ws.Range("E18").Formula = "Female"
$ws.Range("F18").Formula = "2025-11-28 16:28:57"
$ws.Range("G18").Formula = "{
  ""portion"": 0.2,
  ""diet"": 0.42857142857142855,
  ""salt"": 0.4,
  ""fat"": 0.6,
  ""natural"": 0.4,
  ""convenience"": 0.2,
  ""price"": 0.2
}"
$ws.Range("H18").Formula = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I18").Formula = "'0.487"
$ws.Range("J18").Formula = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K18").Formula = "Maruchan Ramen Sabor Pollo"
$ws.Range("L18").Formula = "'0.461"
$ws.Range("M18").Formula = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("N18").Formula = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O18").Formula = "'0.389"
$ws.Range("P18").Formula = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("Q18").Formula = "Kraft Macaroni & Cheese Dinner"
$ws.Range("R18").Formula = "'0.687"
$ws.Range("S18").Formula = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("T18").Formula = "Annie’s Shells & White Cheddar"
$ws.Range("U18").Formula = "'0.649"
$ws.Range("V18").Formula = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("W18").Formula = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("X18").Formula = "'0.582"
$ws.Range("Y18").Formula = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("Z18").Formula = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA18").Formula = "'0.762"
$ws.Range("AB18").Formula = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC18").Formula = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD18").Formula = "'0.635"
$ws.Range("AE18").Formula = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AF18").Formula = "Jack Link’s Beef Jerky Original"
$ws.Range("AG18").Formula = "'0.605"
$ws.Range("AH18").Formula = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# Excel recalcula la altura de la fila al escribir texto multilínea (G18);
# AutoFit la deja en la altura estándar, igual que el resto de filas.
$ws.Rows.Item(18).AutoFit()
